$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7188
$ws.Range("C3").Value = 167974
$ws.Range("C4").Value = 158855
$ws.Range("C8").Value = 65.37
